$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Small wording fixes (plain Find & Replace)
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Statistiker, und Data Scientists", $true, $false, $false, $false, $false,
    $true, 1, $false, "Statistiker und Data Scientists", 2) | Out-Null

$d.Content.Find.Execute(
    "Wodur kann man mit ML Wert schaffen?", $true, $false, $false, $false, $false,
    $true, 1, $false, "Wodurch kann man mit ML Wert schaffen?", 2) | Out-Null

$d.Content.Find.Execute(
    "verschiedenen Algorithmus Kategorien?", $true, $false, $false, $false, $false,
    $true, 1, $false, "verschiedenen Algorithmus-Kategorien?", 2) | Out-Null

$d.Content.Find.Execute(
    "Was sind die Vorteile davon, ein komplexes Input-Output Problem in einfachere Teilprobleme zu zerlegen?",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Was für Vorteile hat es, ein komplexes Input-Output Problem in einfachere Teilprobleme zu zerlegen?", 2) | Out-Null

$d.Content.Find.Execute(
    "Warum könnte ein Modell trotzdem falsch sein, auch wenn es korrekte Vorhersagen für neue Testpunkte liefert?",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Warum kann ein Modell trotzdem falsch sein, auch wenn es korrekte Vorhersagen für Datenpunkte aus dem Testset generiert?", 2) | Out-Null

$d.Content.Find.Execute(
    "Wie kann es passieren, dass ein Modell diskriminiert?", $true, $false, $false, $false, $false,
    $true, 1, $false, "Warum kann es passieren, dass ein Modell diskriminiert und wie erkennt man das?", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) New question block inserted right before the "Fazit" heading:
#      - one Heading 4 paragraph with the new question (+ bookmark)
#      - two blank Normal paragraphs
# ---------------------------------------------------------------------------

$fazitIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($candidate -eq "Fazit") {
        $fazitIndex = $i
        break
    }
}

if ($fazitIndex -ne -1) {
    $anchor = $d.Paragraphs.Item($fazitIndex).Range
    $anchor.Collapse(1)

    $anchor.InsertParagraphBefore()
    $anchor.InsertParagraphBefore()
    $anchor.InsertParagraphBefore()

    $questionPara = $d.Paragraphs.Item($fazitIndex)
    $questionPara.Style = "Heading 4"
    $questionPara.Range.Text = "Aus welchen Gründen könnte es in deiner Domäne / bei deinem nächsten Projekt zum Daten oder Konzept Drift kommen?"

    $blankPara1 = $d.Paragraphs.Item($fazitIndex + 1)
    $blankPara1.Style = "Normal"

    $blankPara2 = $d.Paragraphs.Item($fazitIndex + 2)
    $blankPara2.Style = "Normal"

    $bookmarkRange = $d.Paragraphs.Item($fazitIndex).Range.Duplicate
    $bookmarkRange.Collapse(1)
    $d.Bookmarks.Add("_1bgsj3u87akt", $bookmarkRange) | Out-Null
}
